$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 2 and 3), pushing
# the existing data rows (old 2..42) down to (new 4..44).
$ws.Rows("2:3").Insert()

# The Insert() operation copies formatting from the row above (the bold
# header row). Clear that inherited formatting from the two new rows so
# they look like ordinary data rows, then re-apply the date number format
# used by the rest of column D.
$ws.Rows("2:3").ClearFormats()
$ws.Range("D2:D3").NumberFormat = $ws.Range("D4").NumberFormat

# Populate new row 2
$ws.Cells.Item(2,1).Value2 = 9
$ws.Cells.Item(2,2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(2,3).Value2 = "Metropolitana"
$ws.Cells.Item(2,4).Value2 = 44545
$ws.Cells.Item(2,5).Value2 = 13
$ws.Cells.Item(2,6).Value2 = "Fruta"
$ws.Cells.Item(2,7).Value2 = 100103
$ws.Cells.Item(2,8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(2,9).Value2 = 100103003
$ws.Cells.Item(2,10).Value2 = "Damasco"
$ws.Cells.Item(2,11).Value2 = "Castle Brite"
$ws.Cells.Item(2,12).Value2 = "Especial"
$ws.Cells.Item(2,13).Value2 = 310
$ws.Cells.Item(2,14).Value2 = 18000
$ws.Cells.Item(2,15).Value2 = 18000
$ws.Cells.Item(2,16).Value2 = 18000
$ws.Cells.Item(2,17).Value2 = "`$/caja 18 kilos"
$ws.Cells.Item(2,18).Value2 = "Provincia de Los Andes"
$ws.Cells.Item(2,19).Value2 = 1000
$ws.Cells.Item(2,20).Value2 = 18

# Populate new row 3
$ws.Cells.Item(3,1).Value2 = 9
$ws.Cells.Item(3,2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(3,3).Value2 = "Metropolitana"
$ws.Cells.Item(3,4).Value2 = 44545
$ws.Cells.Item(3,5).Value2 = 13
$ws.Cells.Item(3,6).Value2 = "Fruta"
$ws.Cells.Item(3,7).Value2 = 100103
$ws.Cells.Item(3,8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(3,9).Value2 = 100103003
$ws.Cells.Item(3,10).Value2 = "Damasco"
$ws.Cells.Item(3,11).Value2 = "Castle Brite"
$ws.Cells.Item(3,12).Value2 = "Primera"
$ws.Cells.Item(3,13).Value2 = 350
$ws.Cells.Item(3,14).Value2 = 14400
$ws.Cells.Item(3,15).Value2 = 14400
$ws.Cells.Item(3,16).Value2 = 14400
$ws.Cells.Item(3,17).Value2 = "`$/caja 18 kilos"
$ws.Cells.Item(3,18).Value2 = "Provincia de Los Andes"
$ws.Cells.Item(3,19).Value2 = 800
$ws.Cells.Item(3,20).Value2 = 18

Write-Host "Edit complete"
